$d = $word.ActiveDocument

# --- Edit 1: remove the stray extra space run between
# "Luxury Rental Network " and "system with searching, browsing and" ---
$range1 = $d.Content
$found1 = $range1.Find.Execute("Luxury Rental Network ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $spaceRange = $d.Range($range1.End, $range1.End + 1)
    if ($spaceRange.Text -eq " ") {
        $spaceRange.Delete()
    }
}

# --- Edit 2: replace "Eric Forte took the lead on the writing of this document "
# with "Eric Forte is the owner of this document and, as of this version, is the sole contributor." ---
$range2 = $d.Content
$found2 = $range2.Find.Execute("Eric Forte took the lead on the writing of this document ", $true, $false, $false, $false, $false, $true, 1, $false, "Eric Forte is the owner of this document and, as of this version, is the sole contributor.", 2)

Write-Output "edit1=$found1 edit2=$found2"
